$wb = $excel.ActiveWorkbook

# --- Rename sheets (timestamp refresh) ---
$wb.Worksheets.Item("GNG_TO-16509961231865957").Name = "GNG_TO-1651168722318007"
$wb.Worksheets.Item("NB_TO-16509961253865628").Name = "NB_TO-16511687234797113"
$wb.Worksheets.Item("RS_TO-16509961253865628").Name = "RS_TO-16511687234806914"
$wb.Worksheets.Item("TOL_TO-16509961254585967").Name = "TOL_TO-16511687235428195"
$wb.Worksheets.Item("vSAT_TO-16509961255385635").Name = "vSAT_TO-16511687236053805"

# --- Sheet 1: GNG_TO ---
$ws1 = $wb.Worksheets.Item("GNG_TO-1651168722318007")
$ws1.Range("B2").Value = "go_stims-1651168722276814.csv"
$ws1.Range("B3").Value = "GNG_stims-16511687223018503.csv"
$ws1.Range("B4").Value = "go_stims-1651168722302853.csv"
$ws1.Range("B5").Value = "GNG_stims-16511687223170378.csv"

# --- Sheet 2: NB_TO ---
$ws2 = $wb.Worksheets.Item("NB_TO-16511687234797113")
$ws2.Range("B2").Value = "TB-16511687233618898.csv"
$ws2.Range("B3").Value = "ZB-match_4-16511687225956082.csv"
$ws2.Range("B4").Value = "TB-16511687234641037.csv"
$ws2.Range("B5").Value = "ZB-match_2-1651168722330717.csv"
$ws2.Range("B6").Value = "OB-16511687231213024.csv"
$ws2.Range("B7").Value = "TB-16511687233939748.csv"
$ws2.Range("B8").Value = "OB-165116872280387.csv"
$ws2.Range("B9").Value = "OB-16511687230825295.csv"
$ws2.Range("B10").Value = "ZB-match_1-1651168722472173.csv"

# --- Sheet 3: RS_TO ---
$ws3 = $wb.Worksheets.Item("RS_TO-16511687234806914")
$ws3.Range("B2").Value = "eyes open"
$ws3.Range("B3").Value = "eyes closed"

# --- Sheet 4: TOL_TO ---
$ws4 = $wb.Worksheets.Item("TOL_TO-16511687235428195")
$ws4.Range("B2").Value = "MM_stims-16511687235111825.csv"
$ws4.Range("B3").Value = "ZM_stims-16511687234844882.csv"
$ws4.Range("B4").Value = "MM_stims-16511687235268407.csv"
$ws4.Range("B5").Value = "ZM_stims-16511687235111825.csv"
$ws4.Range("B6").Value = "MM_stims-16511687235428195.csv"
$ws4.Range("B7").Value = "ZM_stims-16511687235278168.csv"

# --- Sheet 5: vSAT_TO ---
$ws5 = $wb.Worksheets.Item("vSAT_TO-16511687236053805")
$ws5.Range("B2").Value = "SAT_stims-165116872355821.csv"
$ws5.Range("B3").Value = "vSAT_stims-16511687235735846.csv"
$ws5.Range("B4").Value = "SAT_stims-16511687235477664.csv"
$ws5.Range("B5").Value = "vSAT_stims-1651168723590659.csv"
